$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2057716666666667
$ws.Range("H2").Value = 0.6173149999999999
$ws.Range("I2").Value = 0.01089677771948535
$ws.Range("J2").Value = 0.01089677771948535
$ws.Range("M2").Value = 6.603177
$ws.Range("N2").Value = 19.809531
$ws.Range("O2").Value = 0.5135477412645301
$ws.Range("P2").Value = 0.5135477412645302
$ws.Range("Q2").Value = 1.358746736585
$ws.Range("R2").Value = 12.228720629265
$ws.Range("S2").Value = 0.005596015584903357
$ws.Range("T2").Value = 0.005596015584903358
$ws.Range("G3").Value = 0.2057716666666667
$ws.Range("H3").Value = 0.6173149999999999
$ws.Range("I3").Value = 0.01089677771948535
$ws.Range("J3").Value = 0.01089677771948535
$ws.Range("O3").Value = 0.02944398858046029
$ws.Range("P3").Value = 0.0294439885804603
$ws.Range("Q3").Value = 0.07790302669277777
$ws.Range("R3").Value = 0.7011272402349999
$ws.Range("S3").Value = 0.0003208445987363407
$ws.Range("T3").Value = 0.0003208445987363407
$ws.Range("G4").Value = 0.2057716666666667
$ws.Range("H4").Value = 0.6173149999999999
$ws.Range("I4").Value = 0.01089677771948535
$ws.Range("J4").Value = 0.01089677771948535
$ws.Range("M4").Value = 3.441487333333333
$ws.Range("N4").Value = 10.324462
$ws.Range("O4").Value = 0.2676541983690312
$ws.Range("P4").Value = 0.2676541983690313
$ws.Range("Q4").Value = 0.7081605843922222
$ws.Range("R4").Value = 6.37344525953
$ws.Range("S4").Value = 0.00291656830531437
$ws.Range("T4").Value = 0.002916568305314371
$ws.Range("G5").Value = 0.2057716666666667
$ws.Range("H5").Value = 0.6173149999999999
$ws.Range("I5").Value = 0.01089677771948535
$ws.Range("J5").Value = 0.01089677771948535
$ws.Range("M5").Value = 2.434707333333333
$ws.Range("N5").Value = 7.304122
$ws.Range("O5").Value = 0.1893540717859783
$ws.Range("P5").Value = 0.1893540717859783
$ws.Range("Q5").Value = 0.5009937858255555
$ws.Range("R5").Value = 4.50894407243
$ws.Range("S5").Value = 0.002063349230531277
$ws.Range("T5").Value = 0.002063349230531278
$ws.Range("I6").Value = 0.01769706320706529
$ws.Range("J6").Value = 0.01769706320706529
$ws.Range("M6").Value = 6.603177
$ws.Range("N6").Value = 19.809531
$ws.Range("O6").Value = 0.5135477412645301
$ws.Range("P6").Value = 0.5135477412645302
$ws.Range("Q6").Value = 2.206691509981
$ws.Range("R6").Value = 19.860223589829
$ws.Range("S6").Value = 0.009088286837003998
$ws.Range("T6").Value = 0.009088286837004004
$ws.Range("I7").Value = 0.01769706320706529
$ws.Range("J7").Value = 0.01769706320706529
$ws.Range("O7").Value = 0.02944398858046029
$ws.Range("P7").Value = 0.0294439885804603
$ws.Range("S7").Value = 0.0005210721269765143
$ws.Range("T7").Value = 0.0005210721269765145
$ws.Range("I8").Value = 0.01769706320706529
$ws.Range("J8").Value = 0.01769706320706529
$ws.Range("M8").Value = 3.441487333333333
$ws.Range("N8").Value = 10.324462
$ws.Range("O8").Value = 0.2676541983690312
$ws.Range("P8").Value = 0.2676541983690313
$ws.Range("Q8").Value = 1.150098033139778
$ws.Range("R8").Value = 10.350882298258
$ws.Range("S8").Value = 0.004736693266173136
$ws.Range("T8").Value = 0.004736693266173138
$ws.Range("I9").Value = 0.01769706320706529
$ws.Range("J9").Value = 0.01769706320706529
$ws.Range("M9").Value = 2.434707333333333
$ws.Range("N9").Value = 7.304122
$ws.Range("O9").Value = 0.1893540717859783
$ws.Range("P9").Value = 0.1893540717859783
$ws.Range("Q9").Value = 0.8136459164664444
$ws.Range("R9").Value = 7.322813248198
$ws.Range("S9").Value = 0.003351010976911636
$ws.Range("T9").Value = 0.003351010976911637
$ws.Range("G10").Value = 0.4895776666666666
$ws.Range("H10").Value = 1.468733
$ws.Range("I10").Value = 0.02592591631545138
$ws.Range("J10").Value = 0.02592591631545138
$ws.Range("M10").Value = 6.603177
$ws.Range("N10").Value = 19.809531
$ws.Range("O10").Value = 0.5135477412645301
$ws.Range("P10").Value = 0.5135477412645302
$ws.Range("Q10").Value = 3.232767988247
$ws.Range("R10").Value = 29.094911894223
$ws.Range("S10").Value = 0.01331419576401328
$ws.Range("T10").Value = 0.01331419576401329
$ws.Range("G11").Value = 0.4895776666666666
$ws.Range("H11").Value = 1.468733
$ws.Range("I11").Value = 0.02592591631545138
$ws.Range("J11").Value = 0.02592591631545138
$ws.Range("O11").Value = 0.02944398858046029
$ws.Range("P11").Value = 0.0294439885804603
$ws.Range("Q11").Value = 0.1853490456307778
$ws.Range("R11").Value = 1.668141410677
$ws.Range("S11").Value = 0.0007633623839301196
$ws.Range("T11").Value = 0.0007633623839301198
$ws.Range("G12").Value = 0.4895776666666666
$ws.Range("H12").Value = 1.468733
$ws.Range("I12").Value = 0.02592591631545138
$ws.Range("J12").Value = 0.02592591631545138
$ws.Range("M12").Value = 3.441487333333333
$ws.Range("N12").Value = 10.324462
$ws.Range("O12").Value = 0.2676541983690312
$ws.Range("P12").Value = 0.2676541983690313
$ws.Range("Q12").Value = 1.684875338516222
$ws.Range("R12").Value = 15.163878046646
$ws.Range("S12").Value = 0.006939180348394727
$ws.Range("T12").Value = 0.006939180348394728
$ws.Range("G13").Value = 0.4895776666666666
$ws.Range("H13").Value = 1.468733
$ws.Range("I13").Value = 0.02592591631545138
$ws.Range("J13").Value = 0.02592591631545138
$ws.Range("M13").Value = 2.434707333333333
$ws.Range("N13").Value = 7.304122
$ws.Range("O13").Value = 0.1893540717859783
$ws.Range("P13").Value = 0.1893540717859783
$ws.Range("Q13").Value = 1.191978335269555
$ws.Range("R13").Value = 10.727805017426
$ws.Range("S13").Value = 0.004909177819113247
$ws.Range("T13").Value = 0.004909177819113248
$ws.Range("G14").Value = 17.85418133333333
$ws.Range("H14").Value = 53.562544
$ws.Range("I14").Value = 0.9454802427579979
$ws.Range("J14").Value = 0.945480242757998
$ws.Range("M14").Value = 6.603177
$ws.Range("N14").Value = 19.809531
$ws.Range("O14").Value = 0.5135477412645301
$ws.Range("P14").Value = 0.5135477412645302
$ws.Range("Q14").Value = 117.894319534096
$ws.Range("R14").Value = 1061.048875806864
$ws.Range("S14").Value = 0.4855492430786094
$ws.Range("T14").Value = 0.4855492430786095
$ws.Range("G15").Value = 17.85418133333333
$ws.Range("H15").Value = 53.562544
$ws.Range("I15").Value = 0.9454802427579979
$ws.Range("J15").Value = 0.945480242757998
$ws.Range("O15").Value = 0.02944398858046029
$ws.Range("P15").Value = 0.0294439885804603
$ws.Range("Q15").Value = 6.759408559592889
$ws.Range("R15").Value = 60.834677036336
$ws.Range("S15").Value = 0.02783870947081732
$ws.Range("T15").Value = 0.02783870947081732
$ws.Range("G16").Value = 17.85418133333333
$ws.Range("H16").Value = 53.562544
$ws.Range("I16").Value = 0.9454802427579979
$ws.Range("J16").Value = 0.945480242757998
$ws.Range("M16").Value = 3.441487333333333
$ws.Range("N16").Value = 10.324462
$ws.Range("O16").Value = 0.2676541983690312
$ws.Range("P16").Value = 0.2676541983690313
$ws.Range("Q16").Value = 61.44493890570311
$ws.Range("R16").Value = 553.004450151328
$ws.Range("S16").Value = 0.253061756449149
$ws.Range("T16").Value = 0.253061756449149
$ws.Range("G17").Value = 17.85418133333333
$ws.Range("H17").Value = 53.562544
$ws.Range("I17").Value = 0.9454802427579979
$ws.Range("J17").Value = 0.945480242757998
$ws.Range("M17").Value = 2.434707333333333
$ws.Range("N17").Value = 7.304122
$ws.Range("O17").Value = 0.1893540717859783
$ws.Range("P17").Value = 0.1893540717859783
$ws.Range("Q17").Value = 43.46970622292978
$ws.Range("R17").Value = 391.227356006368
$ws.Range("S17").Value = 0.1790305337594221
$ws.Range("T17").Value = 0.1790305337594222
